# Update of Excel Modules Files
# ------------------------------------------------------------------
# The sheet originally had columns A:F (Code, Name, ChefModule,
# ElementName1, ElementName2, ElementName3) with only column A filled in
# for rows 2-13 (module codes GSEA11..GSEA26). This edit removes the
# "ElementName3" header (column F becomes unused/empty) and fills in
# the previously-empty columns B:E for every data row with
# member/name/detail text.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "ElementName3" header that used to live in F1. We only clear
# the cell contents (not delete the whole column) so the left-over
# column-width metadata for column F is preserved, matching how Excel
# behaves when you just delete a header's text.
$ws.Range("F1").ClearContents()

# ---- Fill in the new data for columns B (Name), C (ChefModule),
# ---- D (ElementName1) and E (ElementName2) on every existing row.
$ws.Range("B2").Value = 'pede. Suspendisse dui.'
$ws.Range("C2").Value = 'EL Haddad'
$ws.Range("D2").Value = 'Nullam feugiat placerat'
$ws.Range("E2").Value = 'varius et, euismod'

$ws.Range("B3").Value = 'a nunc. In'
$ws.Range("C3").Value = 'Badir'
$ws.Range("D3").Value = 'sodales nisi magna'
$ws.Range("E3").Value = 'elementum sem, vitae'

$ws.Range("B4").Value = 'amet metus. Aliquam'
$ws.Range("C4").Value = 'Ezzine'
$ws.Range("D4").Value = 'Cras vulputate velit'
$ws.Range("E4").Value = 'scelerisque neque sed'

$ws.Range("B5").Value = 'quam vel sapien'
$ws.Range("C5").Value = 'El Alami Hassoun'
$ws.Range("D5").Value = 'Nunc mauris elit,'
$ws.Range("E5").Value = 'libero et tristique'

$ws.Range("B6").Value = 'feugiat nec, diam.'
$ws.Range("C6").Value = 'Lazaar'
$ws.Range("D6").Value = 'pellentesque. Sed dictum.'
$ws.Range("E6").Value = 'ridiculus mus. Proin'

$ws.Range("B7").Value = 'nonummy. Fusce fermentum'
$ws.Range("C7").Value = 'El Haddad'
$ws.Range("D7").Value = 'neque pellentesque massa'
$ws.Range("E7").Value = 'Mauris eu turpis.'

$ws.Range("B8").Value = 'a, arcu. Sed'
$ws.Range("C8").Value = 'EL Haddad'
$ws.Range("D8").Value = 'sit amet risus.'
$ws.Range("E8").Value = 'Nulla facilisi. Sed'

$ws.Range("B9").Value = 'Suspendisse eleifend. Cras'
$ws.Range("C9").Value = 'El Alami Hassoun'
$ws.Range("D9").Value = 'velit dui, semper'
$ws.Range("E9").Value = 'ligula elit, pretium'

$ws.Range("B10").Value = 'ante. Nunc mauris'
$ws.Range("C10").Value = 'Badir'
$ws.Range("D10").Value = 'tortor at risus.'
$ws.Range("E10").Value = 'felis. Donec tempor,'

$ws.Range("B11").Value = 'lobortis quam a'
$ws.Range("C11").Value = 'Ezzine'
$ws.Range("D11").Value = 'euismod est arcu'
$ws.Range("E11").Value = 'ligula eu enim.'

$ws.Range("B12").Value = 'rhoncus. Nullam velit'
$ws.Range("C12").Value = 'Ben Achrab'
$ws.Range("D12").Value = 'ut dolor dapibus'
$ws.Range("E12").Value = 'commodo tincidunt nibh.'

$ws.Range("B13").Value = 'Donec tincidunt. Donec'
$ws.Range("C13").Value = 'EL Haddad'
$ws.Range("D13").Value = 'ornare tortor at'
$ws.Range("E13").Value = 'ac, feugiat non,'

# ---- B2 carries an explicit (non-theme) black font, same as the
# ---- rest of the author's formatting pass.
$ws.Range("B2").Font.ThemeFont = 0
$ws.Range("B2").Font.Color = 0

# ---- Column widths: B:E were auto-fit to the new, wider content;
# ---- column F keeps its old width (it's simply no longer used).
$ws.Columns.Item(2).ColumnWidth = 26.166666666666668
$ws.Columns.Item(3).ColumnWidth = 15.307291666666666
$ws.Columns.Item(4).ColumnWidth = 24.307291666666668
$ws.Columns.Item(5).ColumnWidth = 22.736979166666668

# ---- Final selection left on E9 (matches the saved view state).
$ws.Range("E9").Select()
